# Update transition-probability matrix on Sheet1 ("Spalding_B") with
# recomputed values after adding more simulated games (commit: "added more
# games, sped up simulate game logic, and drafted optimization logic").
# Row labels (col A) are state-from codes (Af0..Br0); column headers (row 1)
# are state-to indices 0..17 mapped to columns B..S.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Af0)
$ws.Range("B2").Value = 0.1481481481481481
$ws.Range("C2").Value = 0.5925925925925926
$ws.Range("J2").Value = 0.03703703703703703
$ws.Range("O2").Value = 0.03703703703703703
$ws.Range("P2").Value = 0.1111111111111111
$ws.Range("S2").Value = 0.07407407407407407

# Row 3 (Af1)
$ws.Range("C3").Value = 0.05882352941176471
$ws.Range("J3").Value = 0.1176470588235294
$ws.Range("P3").Value = 0.7058823529411765
$ws.Range("S3").Value = 0.1176470588235294

# Row 4 (Af2)
$ws.Range("J4").Value = 0.1428571428571428
$ws.Range("P4").Value = 0.5714285714285714
$ws.Range("S4").Value = 0.2857142857142857

# Row 6 (Ai0)
$ws.Range("B6").Value = 0.08333333333333333
$ws.Range("F6").Value = 0.08333333333333333
$ws.Range("J6").Value = 0.4166666666666667
$ws.Range("O6").Value = 0.08333333333333333
$ws.Range("Q6").Value = 0.08333333333333333
$ws.Range("S6").Value = 0.25

# Row 7 (Ai1)
$ws.Range("B7").Value = 0.4166666666666667
$ws.Range("D7").Value = 0.08333333333333333
$ws.Range("J7").Value = 0.08333333333333333
$ws.Range("Q7").Value = 0.1666666666666667
$ws.Range("R7").Value = 0.08333333333333333
$ws.Range("S7").Value = 0.1666666666666667

# Row 8 (Ai2)
$ws.Range("B8").Value = 0.03125
$ws.Range("D8").Value = 0.0625
$ws.Range("F8").Value = 0.0625
$ws.Range("J8").Value = 0.28125
$ws.Range("O8").Value = 0.0625
$ws.Range("Q8").Value = 0.15625
$ws.Range("R8").Value = 0.125
$ws.Range("S8").Value = 0.21875

# Row 9 (Ai3)
$ws.Range("D9").Value = 0.1666666666666667
$ws.Range("J9").Value = 0.1666666666666667
$ws.Range("Q9").Value = 0.1666666666666667

# Row 10 (Ar0)
$ws.Range("B10").Value = 0.1588785046728972
$ws.Range("D10").Value = 0.02803738317757009
$ws.Range("F10").Value = 0.02803738317757009
$ws.Range("J10").Value = 0.1308411214953271
$ws.Range("O10").Value = 0.03738317757009346
$ws.Range("Q10").Value = 0.205607476635514
$ws.Range("R10").Value = 0.04672897196261682
$ws.Range("S10").Value = 0.3644859813084112

# Row 11 (Bf0)
$ws.Range("G11").Value = 0.2142857142857143
$ws.Range("J11").Value = 0.07142857142857142

# Row 13 (Bf2)
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.5

# Row 15 (Bi0)
$ws.Range("F15").Value = 0.03846153846153846
$ws.Range("H15").Value = 0.1153846153846154
$ws.Range("J15").Value = 0.4230769230769231
$ws.Range("M15").Value = 0.03846153846153846

# Row 16 (Bi1)
$ws.Range("F16").Value = 0.05263157894736842
$ws.Range("H16").Value = 0.1578947368421053
$ws.Range("J16").Value = 0.5789473684210527
$ws.Range("K16").Value = 0.05263157894736842
$ws.Range("O16").Value = 0.05263157894736842
$ws.Range("S16").Value = 0.1052631578947368

# Row 17 (Bi2)
$ws.Range("F17").Value = 0.0625
$ws.Range("H17").Value = 0.15625
$ws.Range("I17").Value = 0.09375
$ws.Range("J17").Value = 0.53125
$ws.Range("K17").Value = 0.03125
$ws.Range("O17").Value = 0.0625
$ws.Range("S17").Value = 0.0625

# Row 18 (Bi3)
$ws.Range("H18").Value = 0.1
$ws.Range("J18").Value = 0.6
$ws.Range("K18").Value = 0.1
$ws.Range("O18").Value = 0.1
$ws.Range("S18").Value = 0.1

# Row 19 (Br0)
$ws.Range("F19").Value = 0.01282051282051282
$ws.Range("H19").Value = 0.217948717948718
$ws.Range("I19").Value = 0.03846153846153846
$ws.Range("J19").Value = 0.3974358974358974
$ws.Range("K19").Value = 0.1153846153846154
$ws.Range("M19").Value = 0.01282051282051282
$ws.Range("O19").Value = 0.1153846153846974
$ws.Range("S19").Value = 0.08974358974358974
